$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912" : new scrape row + updated totals ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:19:33"
$ws1.Range("A3").Value = "Total filas: 2"

$ws1.Range("A6").Value = "01:19:33"
$ws1.Range("D6").Value = 6

$ws1.Range("A7").Value = "01:19:33"
$ws1.Range("B7").Value = "03:01"
$ws1.Range("C7").Value = "15_ABASTO"
$ws1.Range("D7").Value = 102
$ws1.Range("E7").Value = "LP1912"

# --- Sheet "LP1912-215" : refreshed timestamp only ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 01:19:33"

# --- Sheet "6203-6173" : refreshed timestamp only ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 01:19:33"
